# Added Profile Description Test cases - on progress
#
# Inserts a new "Description" column (G) on the Profile sheet, in between
# the existing "EarnTarget" column (F) and the "Country" column (old G,
# now H), and fills in a handful of description test-case values used for
# profile-description validation testing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Profile")

# Shift the old G:O columns right to make room for the new Description
# column; the new column inherits formatting from its left neighbour (F).
$ws.Columns("G:G").Insert()

# Header
$ws.Range("G1").Value = "Description"
$ws.Range("G1").Interior.ColorIndex = $ws.Range("H1").Interior.ColorIndex

# Test-case values
$ws.Range("G2").Value = "I love coding and working on my skill to improve and get better to become a Software Tester."

$longDescription = @'
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
This is another test to check if description will accept more than 598 characters.
A total of 597 characters including spaces and periods.
'@

$ws.Range("G3").Value = $longDescription
$ws.Range("G4").Value = "Invalid characters"
$ws.Range("G5").Value = "Empty description"
$ws.Range("G6").Value = "description starts with a space followed by a valid description"

# Widen the Description column so the longest entry is visible, mirroring
# the bestFit auto-size Excel applies when a column is double-clicked.
$ws.Columns("G:G").ColumnWidth = 255.6328125

$ws.Range("G9").Select()
